# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gains a new (blank) column at position N
# (the 14th column), pushing the existing "Late" / "heading" / "Outstanding"
# columns one place to the right (N->O, O->P, P->Q). The active sheet/tab
# also moves from "Transactions" back to "Repayment schedule", with the
# selection left on the new last column of row 5 (Q5).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column). Excel shifts the
# old N/O/P columns (and their styles/values) one place to the right
# automatically.
$ws.Columns.Item(14).Insert()

# Match the width of the freshly inserted column to the one used for the
# rest of the repayment-schedule table (same width class as column M).
$ws.Columns.Item(14).ColumnWidth = 9.83

# Re-select the "Repayment schedule" tab (it was "Transactions" before) and
# leave the selection on Q5, mirroring the saved view state.
$ws.Activate()
$ws.Range("Q5").Select()
